$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "[-, -, -, 'MEC-3B-M. Motor Endot.']"
$ws.Range("D3").Value = "-"

$ws.Range("C4").Value = "[-, -, -, 'MEC-3B-M. Motor Endot.']"
$ws.Range("D4").Value = "-"

$ws.Range("C6").Value = "[-, -, -, 'MEC-3B-M. Motor Endot.']"
$ws.Range("D6").Value = "-"

$ws.Range("C7").Value = "[-, -, -, 'MEC-3B-M. Motor Endot.']"
$ws.Range("D7").Value = "-"

$ws.Range("D10").Value = "-"
$ws.Range("F10").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"

$ws.Range("B11").Value = "-"

$ws.Range("B12").Value = "-"

$ws.Range("D14").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"

$ws.Range("D15").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"

$ws.Range("D16").Value = "-"
$ws.Range("F16").Value = "[-, 'MEC-3A-M. Motor Endot.', -, -]"

$ws.Range("B18").Value = "['MEC-1NA-Manut. Mot. End.', -, -, 'MEC-1NB-Manut. Mot. End.']"
$ws.Range("E18").Value = "-"

$ws.Range("B19").Value = "['MEC-1NA-Manut. Mot. End.', -, -, 'MEC-1NB-Manut. Mot. End.']"
$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "-"

$ws.Range("B20").Value = "[-, -, -, 'MEC-1NB-Manut. Mot. End.']"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "[-, -, 'MEC-1NA-Manut. Mot. End.', -]"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "[-, -, -, 'MEC-1NB-Manut. Mot. End.']"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, -, 'MEC-1NA-Manut. Mot. End.', -]"
$ws.Range("E21").Value = "-"
